$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1").Value = " lower size"
$ws.Range("F2").Select()
